# Add season record columns (Wins, Losses, Ties) to the MIA 2011 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new labels, formatted like the neighboring header cell.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats - reuse AC1's style

# Data rows 2-46: every player row gets the team's season record.
$ws.Range("AD2:AD46").Value = 72
$ws.Range("AE2:AE46").Value = 90
$ws.Range("AF2:AF46").Value = 0
